$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format before writing, to preserve values
# like "1.001" / "317.35" as literal text instead of being auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '28.351.68'
$ws.Range("E2").Value = '  +1.47%  '
$ws.Range("D3").Value = '1.823.06'
$ws.Range("E3").Value = '  +2.49%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '317.35'
$ws.Range("E5").Value = '  +0.46%  '
$ws.Range("E6").Value = '  +0.05%  '
$ws.Range("D7").Value = '0.5335'
$ws.Range("E7").Value = '  -0.97%  '
$ws.Range("D8").Value = '0.4033'
$ws.Range("E8").Value = '  +7.15%  '
$ws.Range("D9").Value = '0.07609'
$ws.Range("E9").Value = '  +2.41%  '
$ws.Range("D10").Value = '41.81'
$ws.Range("E10").Value = '  +0.46%  '
$ws.Range("D11").Value = '1.107'
$ws.Range("E11").Value = '  +1.43%  '
$ws.Range("D12").Value = '6.318'
$ws.Range("E12").Value = '  +4.11%  '
$ws.Range("E13").Value = '  +0.08%  '
$ws.Range("D14").Value = '7.602'
$ws.Range("E14").Value = '  +5.61%  '
$ws.Range("D15").Value = '20.88'
$ws.Range("E15").Value = '  +2.25%  '
$ws.Range("D16").Value = '1.830.36'
$ws.Range("E16").Value = '  +3.34%  '
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").Value = '0.00001075'
$ws.Range("E17").Value = '  +2.29%  '
$ws.Range("B18").Value = 'Litecoin'
$ws.Range("C18").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D18").Value = '89.37'
$ws.Range("E18").Value = '  +1.49%  '
$ws.Range("D19").Value = '0.06595'
$ws.Range("E19").Value = '  +2.47%  '
$ws.Range("D20").Value = '17.67'
$ws.Range("E20").Value = '  +2.73%  '
$ws.Range("D21").Value = '1.001'
$ws.Range("E21").Value = '  +0.06%  '
$ws.Range("D22").Value = '6.084'
$ws.Range("E22").Value = '  +3.68%  '
$ws.Range("D23").Value = '28.373.03'
$ws.Range("E23").Value = '  +1.42%  '
$ws.Range("D24").Value = '11.18'
$ws.Range("E24").Value = '  +0.38%  '
$ws.Range("D25").Value = '2.206'
$ws.Range("E25").Value = '  +5.67%  '
$ws.Range("D26").Value = '2.459'
$ws.Range("E26").Value = '  +7.96%  '
$ws.Range("D27").Value = '157.50'
$ws.Range("E27").Value = '  +0.94%  '
$ws.Range("E28").Value = '  +2.06%  '
$ws.Range("D29").Value = '2.040.65'
$ws.Range("E29").Value = '  +3.03%  '
$ws.Range("D30").Value = '123.89'
$ws.Range("E30").Value = '  +3.45%  '
$ws.Range("D31").Value = '1.122'
$ws.Range("E31").Value = '  +1.49%  '
$ws.Range("D32").Value = '0.1101'
$ws.Range("E32").Value = '  +5.29%  '
$ws.Range("D33").Value = '5.657'
$ws.Range("E33").Value = '  +2.84%  '
$ws.Range("D34").Value = '0.07400'
$ws.Range("E34").Value = '  +15.49%  '
$ws.Range("D35").Value = '3.643'
$ws.Range("E35").Value = '  +0.07%  '
$ws.Range("D36").Value = '0.2233'
$ws.Range("E36").Value = '  -0.91%  '
$ws.Range("D37").Value = '0.02342'
$ws.Range("E37").Value = '  +3.47%  '
$ws.Range("D38").Value = '8.906'
$ws.Range("E38").Value = '  +6.17%  '
$ws.Range("D39").Value = '5.198'
$ws.Range("E39").Value = '  +4.80%  '
$ws.Range("E40").Value = '  +2.29%  '
$ws.Range("D41").Value = '0.6255'
$ws.Range("E41").Value = '  +2.13%  '
$ws.Range("D42").Value = '1.182'
$ws.Range("E42").Value = '  +0.66%  '
$ws.Range("E43").Value = '  +0.06%  '
$ws.Range("D44").Value = '1.394'
$ws.Range("E44").Value = '  -2.35%  '
$ws.Range("D45").Value = '13.52'
$ws.Range("E45").Value = '  +1.80%  '
$ws.Range("D46").Value = '3.700'
$ws.Range("E46").Value = '  +1.13%  '
$ws.Range("D47").Value = '0.5836'
$ws.Range("E47").Value = '  +1.87%  '
$ws.Range("D48").Value = '125.05'
$ws.Range("E48").Value = '  -1.00%  '
$ws.Range("D49").Value = '1.989'
$ws.Range("E49").Value = '  +3.65%  '
$ws.Range("D50").Value = '1.198'
$ws.Range("E50").Value = '  +1.72%  '
$ws.Range("D51").Value = '0.06895'
$ws.Range("E51").Value = '  +1.66%  '

# Restore column D cell style to the workbook default (removes the explicit
# Text number format so the saved style matches the original, unstyled cells).
$ws.Range("D2:D51").Style = "Normal"
